# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos list
# with refreshed figures, and swaps the ImmutableX / PancakeSwap rows
# (22 <-> 23), matching the GitHub Actions crypto-data refresh commit.
#
# Note: several Price values (e.g. "312.85") look like plain numbers to
# Excel's automatic type detection, which would silently convert the cell
# from Text to Number on assignment (losing the trailing zero / fixed
# formatting the source data relies on). Those are written with a leading
# apostrophe - exactly like a user typing '312.85 into the cell - to force
# a Text value, then the transient "quote prefix" style that trick leaves
# behind is cleared via Style = 'Normal' so the cell's formatting is left
# exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.238.12'
$ws.Range('E2').Value = '  +2.11%  '
$ws.Range('D3').Value = '2.348.78'
$ws.Range('E3').Value = '  +6.12%  '
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').Value = "'312.85"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.91%  '
$ws.Range('D6').Value = "'109.55"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.13%  '
$ws.Range('E7').Value = '  +3.49%  '
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').Value = "'0.635"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.87%  '
$ws.Range('D10').Value = "'42.99"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.04%  '
$ws.Range('D11').Value = "'0.0938"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.34%  '
$ws.Range('D12').Value = "'8.87"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.46%  '
$ws.Range('E13').Value = '  +9.43%  '
$ws.Range('E14').Value = '  +2.30%  '
$ws.Range('D15').Value = "'16.30"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +9.29%  '
$ws.Range('D16').Value = '2.704.67'
$ws.Range('E16').Value = '  +6.13%  '
$ws.Range('D17').Value = '2.346.31'
$ws.Range('E17').Value = '  +4.32%  '
$ws.Range('D18').Value = '43.239.70'
$ws.Range('E19').Value = '  +3.27%  '
$ws.Range('E20').Value = '  -1.25%  '
$ws.Range('D21').Value = "'75.44"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.94%  '
$ws.Range('B22').Value = 'ImmutableX'
$ws.Range('C22').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D22').Value = "'2.59"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +13.51%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').Value = "'3.45"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.42%  '
$ws.Range('D24').Value = "'255.21"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +12.24%  '
$ws.Range('D25').Value = "'9.11"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.21%  '
$ws.Range('D26').Value = "'12.04"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.56%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').Value = '  +2.24%  '
$ws.Range('D29').Value = "'2.26"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('D30').Value = "'22.38"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +7.42%  '
$ws.Range('D31').Value = "'174.49"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('E32').Value = '  -1.05%  '
$ws.Range('D33').Value = "'0.0926"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.57%  '
$ws.Range('E34').Value = '  +10.35%  '
$ws.Range('E35').Value = '  +6.27%  '
$ws.Range('D36').Value = "'4.96"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.29%  '
$ws.Range('E37').Value = '  +3.89%  '
$ws.Range('E38').Value = '  -3.52%  '
$ws.Range('E39').Value = '  +2.66%  '
$ws.Range('D40').Value = "'2.70"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.54%  '
$ws.Range('D41').Value = "'72.72"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.09%  '
$ws.Range('E42').Value = '  +14.94%  '
$ws.Range('E43').Value = '  +1.97%  '
$ws.Range('D44').Value = "'12.77"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.00%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  +4.14%  '
$ws.Range('E47').Value = '  +12.25%  '
$ws.Range('D48').Value = "'111.02"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.61%  '
$ws.Range('E49').Value = '  -0.63%  '
$ws.Range('E50').Value = '  +4.18%  '
$ws.Range('D51').Value = "'69.86"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.56%  '
